$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1 - copy format from E1 (bold/centered/bordered header style)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value = "time_taken"

# time_taken values for rows 2-48 (default style, plain text)
$ws.Cells.Item(2, 6).Value = "2021-10-05 10:52:10.775495"
$ws.Cells.Item(3, 6).Value = "2021-10-05 10:52:10.775506"
$ws.Cells.Item(4, 6).Value = "2021-10-05 10:52:10.775510"
$ws.Cells.Item(5, 6).Value = "2021-10-05 10:52:10.775512"
$ws.Cells.Item(6, 6).Value = "2021-10-05 10:52:10.775515"
$ws.Cells.Item(7, 6).Value = "2021-10-05 10:52:10.775518"
$ws.Cells.Item(8, 6).Value = "2021-10-05 10:52:10.775520"
$ws.Cells.Item(9, 6).Value = "2021-10-05 10:52:10.775523"
$ws.Cells.Item(10, 6).Value = "2021-10-05 10:52:10.775526"
$ws.Cells.Item(11, 6).Value = "2021-10-05 10:52:10.775528"
$ws.Cells.Item(12, 6).Value = "2021-10-05 10:52:10.775531"
$ws.Cells.Item(13, 6).Value = "2021-10-05 10:52:10.775533"
$ws.Cells.Item(14, 6).Value = "2021-10-05 10:52:10.775536"
$ws.Cells.Item(15, 6).Value = "2021-10-05 10:52:10.775538"
$ws.Cells.Item(16, 6).Value = "2021-10-05 10:52:10.775540"
$ws.Cells.Item(17, 6).Value = "2021-10-05 10:52:10.775543"
$ws.Cells.Item(18, 6).Value = "2021-10-05 10:52:10.775546"
$ws.Cells.Item(19, 6).Value = "2021-10-05 10:52:10.775548"
$ws.Cells.Item(20, 6).Value = "2021-10-05 10:52:10.775551"
$ws.Cells.Item(21, 6).Value = "2021-10-05 10:52:10.775553"
$ws.Cells.Item(22, 6).Value = "2021-10-05 10:52:10.775556"
$ws.Cells.Item(23, 6).Value = "2021-10-05 10:52:10.775558"
$ws.Cells.Item(24, 6).Value = "2021-10-05 10:52:10.775561"
$ws.Cells.Item(25, 6).Value = "2021-10-05 10:52:10.775563"
$ws.Cells.Item(26, 6).Value = "2021-10-05 10:52:10.775566"
$ws.Cells.Item(27, 6).Value = "2021-10-05 10:52:10.775568"
$ws.Cells.Item(28, 6).Value = "2021-10-05 10:52:10.775571"
$ws.Cells.Item(29, 6).Value = "2021-10-05 10:52:10.775573"
$ws.Cells.Item(30, 6).Value = "2021-10-05 10:52:10.775576"
$ws.Cells.Item(31, 6).Value = "2021-10-05 10:52:10.775578"
$ws.Cells.Item(32, 6).Value = "2021-10-05 10:52:10.775581"
$ws.Cells.Item(33, 6).Value = "2021-10-05 10:52:10.775583"
$ws.Cells.Item(34, 6).Value = "2021-10-05 10:52:10.775586"
$ws.Cells.Item(35, 6).Value = "2021-10-05 10:52:10.775589"
$ws.Cells.Item(36, 6).Value = "2021-10-05 10:52:10.775592"
$ws.Cells.Item(37, 6).Value = "2021-10-05 10:52:10.775594"
$ws.Cells.Item(38, 6).Value = "2021-10-05 10:52:10.775597"
$ws.Cells.Item(39, 6).Value = "2021-10-05 10:52:10.775599"
$ws.Cells.Item(40, 6).Value = "2021-10-05 10:52:10.775602"
$ws.Cells.Item(41, 6).Value = "2021-10-05 10:52:10.775604"
$ws.Cells.Item(42, 6).Value = "2021-10-05 10:52:10.775607"
$ws.Cells.Item(43, 6).Value = "2021-10-05 10:52:10.775610"
$ws.Cells.Item(44, 6).Value = "2021-10-05 10:52:10.775612"
$ws.Cells.Item(45, 6).Value = "2021-10-05 10:52:10.775615"
$ws.Cells.Item(46, 6).Value = "2021-10-05 10:52:10.775617"
$ws.Cells.Item(47, 6).Value = "2021-10-05 10:52:10.775620"
$ws.Cells.Item(48, 6).Value = "2021-10-05 10:52:10.775622"
